# Auto-generated Excel COM-interop script
# Applies the "Update countries & provincias Spain" edit to the Pais sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 11:12"

# Row 9
$ws.Cells.Item(9, 1).Value = "Iran"
$ws.Cells.Item(9, 2).Value = 32332
$ws.Cells.Item(9, 3).Value = 2926
$ws.Cells.Item(9, 4).Value = 11133
$ws.Cells.Item(9, 5).Value = 18821
$ws.Cells.Item(9, 6).Value = 2893
$ws.Cells.Item(9, 7).Value = 144
$ws.Cells.Item(9, 8).Value = 2378

# Row 15
$ws.Cells.Item(15, 1).Value = "Belgica"
$ws.Cells.Item(15, 2).Value = 7284
$ws.Cells.Item(15, 3).Value = 1049
$ws.Cells.Item(15, 4).Value = 675
$ws.Cells.Item(15, 5).Value = 6320
$ws.Cells.Item(15, 6).Value = 605
$ws.Cells.Item(15, 7).Value = 69
$ws.Cells.Item(15, 8).Value = 289

# Row 16
$ws.Cells.Item(16, 1).Value = "Austria"
$ws.Cells.Item(16, 2).Value = 7196
$ws.Cells.Item(16, 3).Value = 287
$ws.Cells.Item(16, 4).Value = 225
$ws.Cells.Item(16, 5).Value = 6913
$ws.Cells.Item(16, 6).Value = 96
$ws.Cells.Item(16, 7).Value = 9
$ws.Cells.Item(16, 8).Value = 58

# Row 66
$ws.Cells.Item(66, 1).Value = "Lituania"
$ws.Cells.Item(66, 2).Value = 345
$ws.Cells.Item(66, 3).Value = 46
$ws.Cells.Item(66, 4).Value = 1
$ws.Cells.Item(66, 5).Value = 339
$ws.Cells.Item(66, 6).Value = 2
$ws.Cells.Item(66, 7).Value = 1
$ws.Cells.Item(66, 8).Value = 5

# Row 76
$ws.Cells.Item(76, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(76, 2).Value = 230
$ws.Cells.Item(76, 3).Value = 39
$ws.Cells.Item(76, 4).Value = 5
$ws.Cells.Item(76, 5).Value = 222
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 3

# Row 77
$ws.Cells.Item(77, 1).Value = "Tunez"
$ws.Cells.Item(77, 2).Value = 227
$ws.Cells.Item(77, 3).Value = 30
$ws.Cells.Item(77, 4).Value = 2
$ws.Cells.Item(77, 5).Value = 219
$ws.Cells.Item(77, 6).Value = 10
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 6

# Row 78
$ws.Cells.Item(78, 1).Value = "Eslovaquia"
$ws.Cells.Item(78, 2).Value = 226
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 4).Value = 2
$ws.Cells.Item(78, 5).Value = 224
$ws.Cells.Item(78, 6).Value = 2
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 0

# Row 79
$ws.Cells.Item(79, 1).Value = "Kuwait"
$ws.Cells.Item(79, 2).Value = 225
$ws.Cells.Item(79, 3).Value = 17
$ws.Cells.Item(79, 4).Value = 57
$ws.Cells.Item(79, 5).Value = 168
$ws.Cells.Item(79, 6).Value = 11
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 0

# Row 80
$ws.Cells.Item(80, 1).Value = "Principado de Andorra"
$ws.Cells.Item(80, 2).Value = 224
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = 1
$ws.Cells.Item(80, 5).Value = 220
$ws.Cells.Item(80, 6).Value = 6
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 3

# Row 81
$ws.Cells.Item(81, 1).Value = "Ucrania"
$ws.Cells.Item(81, 2).Value = 218
$ws.Cells.Item(81, 3).Value = 22
$ws.Cells.Item(81, 4).Value = 4
$ws.Cells.Item(81, 5).Value = 209
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 0
$ws.Cells.Item(81, 8).Value = 5

# Row 82
$ws.Cells.Item(82, 1).Value = "Jordania"
$ws.Cells.Item(82, 2).Value = 212
$ws.Cells.Item(82, 3).Value = 0
$ws.Cells.Item(82, 4).Value = 2
$ws.Cells.Item(82, 5).Value = 210
$ws.Cells.Item(82, 6).Value = 0
$ws.Cells.Item(82, 7).Value = 0
$ws.Cells.Item(82, 8).Value = 0

# Row 83
$ws.Cells.Item(83, 1).Value = "San Marino"
$ws.Cells.Item(83, 2).Value = 208
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 4
$ws.Cells.Item(83, 5).Value = 183
$ws.Cells.Item(83, 6).Value = 12
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 21

# Row 90
$ws.Cells.Item(90, 1).Value = "Islas Feroe"
$ws.Cells.Item(90, 2).Value = 144
$ws.Cells.Item(90, 3).Value = 4
$ws.Cells.Item(90, 4).Value = 54
$ws.Cells.Item(90, 5).Value = 90
$ws.Cells.Item(90, 6).Value = 2
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 0

# Row 91
$ws.Cells.Item(91, 1).Value = "Ghana"
$ws.Cells.Item(91, 2).Value = 136
$ws.Cells.Item(91, 3).Value = 4
$ws.Cells.Item(91, 4).Value = 1
$ws.Cells.Item(91, 5).Value = 131
$ws.Cells.Item(91, 6).Value = 1
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 4

# Row 92
$ws.Cells.Item(92, 1).Value = "Reunion"
$ws.Cells.Item(92, 2).Value = 135
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 1
$ws.Cells.Item(92, 5).Value = 134
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0

# Row 93
$ws.Cells.Item(93, 1).Value = "Malta"
$ws.Cells.Item(93, 2).Value = 134
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 2
$ws.Cells.Item(93, 5).Value = 132
$ws.Cells.Item(93, 6).Value = 1
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0

# Row 97
$ws.Cells.Item(97, 1).Value = "Brunei"
$ws.Cells.Item(97, 2).Value = 115
$ws.Cells.Item(97, 3).Value = 1
$ws.Cells.Item(97, 4).Value = 11
$ws.Cells.Item(97, 5).Value = 104
$ws.Cells.Item(97, 6).Value = 1
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0

# Row 104
$ws.Cells.Item(104, 1).Value = "Bielorrusia"
$ws.Cells.Item(104, 2).Value = 94
$ws.Cells.Item(104, 3).Value = 8
$ws.Cells.Item(104, 4).Value = 29
$ws.Cells.Item(104, 5).Value = 65
$ws.Cells.Item(104, 6).Value = 2
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 0

# Row 105
$ws.Cells.Item(105, 1).Value = "Estado de Palestina"
$ws.Cells.Item(105, 2).Value = 91
$ws.Cells.Item(105, 3).Value = 5
$ws.Cells.Item(105, 4).Value = 17
$ws.Cells.Item(105, 5).Value = 73
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 1

# Row 106
$ws.Cells.Item(106, 1).Value = "Camerun"
$ws.Cells.Item(106, 2).Value = 88
$ws.Cells.Item(106, 3).Value = 13
$ws.Cells.Item(106, 4).Value = 2
$ws.Cells.Item(106, 5).Value = 84
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 1
$ws.Cells.Item(106, 8).Value = 2

# Row 122
$ws.Cells.Item(122, 1).Value = "Mayotte"
$ws.Cells.Item(122, 2).Value = 50
$ws.Cells.Item(122, 3).Value = 14
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 50
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 0

# Row 123
$ws.Cells.Item(123, 1).Value = "Ruanda"
$ws.Cells.Item(123, 2).Value = 50
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 50
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 0

# Row 124
$ws.Cells.Item(124, 1).Value = "Banglades"
$ws.Cells.Item(124, 2).Value = 48
$ws.Cells.Item(124, 3).Value = 4
$ws.Cells.Item(124, 4).Value = 11
$ws.Cells.Item(124, 5).Value = 32
$ws.Cells.Item(124, 6).Value = 1
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 5

# Row 125
$ws.Cells.Item(125, 1).Value = "Puerto Rico"
$ws.Cells.Item(125, 2).Value = 39
$ws.Cells.Item(125, 3).Value = 0
$ws.Cells.Item(125, 4).Value = 1
$ws.Cells.Item(125, 5).Value = 36
$ws.Cells.Item(125, 6).Value = 0
$ws.Cells.Item(125, 7).Value = 0
$ws.Cells.Item(125, 8).Value = 2
